$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - force text format so numeric-looking
# strings (e.g. "0.9999", "331.90") are stored as text, matching the
# original inline-string cell type, not auto-converted to numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.137.01"
$ws.Range("D3").Value = "1.923.08"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D5").Value = "331.90"
$ws.Range("D7").Value = "0.5258"
$ws.Range("D8").Value = "0.4064"
$ws.Range("D9").Value = "0.08554"
$ws.Range("D10").Value = "43.31"
$ws.Range("D11").Value = "1.131"
$ws.Range("D12").Value = "22.65"
$ws.Range("D13").Value = "6.436"
$ws.Range("D14").Value = "1.919.76"
$ws.Range("D15").Value = "7.418"
$ws.Range("D16").Value = "0.9998"
$ws.Range("D17").Value = "96.72"
$ws.Range("D18").Value = "0.00001118"
$ws.Range("D19").Value = "0.06719"
$ws.Range("D20").Value = "18.33"
$ws.Range("D21").Value = "0.9994"
$ws.Range("D23").Value = "30.130.74"
$ws.Range("D24").Value = "11.31"
$ws.Range("D26").Value = "2.142.75"
$ws.Range("D27").Value = "21.22"
$ws.Range("D28").Value = "160.34"
$ws.Range("D29").Value = "2.480"
$ws.Range("D30").Value = "129.68"
$ws.Range("D31").Value = "1.084"
$ws.Range("D33").Value = "6.138"
$ws.Range("D34").Value = "3.646"
$ws.Range("D35").Value = "0.02527"
$ws.Range("D36").Value = "0.06610"
$ws.Range("D37").Value = "0.2231"
$ws.Range("D38").Value = "9.107"
$ws.Range("D39").Value = "1.240"
$ws.Range("D40").Value = "5.222"
$ws.Range("D41").Value = "0.6584"
$ws.Range("D42").Value = "11.73"
$ws.Range("D43").Value = "1.245"
$ws.Range("D44").Value = "0.6210"
$ws.Range("D45").Value = "13.30"
$ws.Range("D46").Value = "3.789"
$ws.Range("D47").Value = "2.097"
$ws.Range("D48").Value = "1.250"
$ws.Range("D49").Value = "125.32"
$ws.Range("D50").Value = "80.10"

# Restore default (General) styling so no stray number format is left
# attached to the cell once the text value is safely stored.
foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}

# Volume(1h) (column E) updates - values already contain "%" and padding
# spaces, so Excel stores them as text without any extra handling.
$ws.Range("E2").Value = "  +5.73%  "
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("E5").Value = "  +5.05%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  +3.29%  "
$ws.Range("E8").Value = "  +4.14%  "
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").Value = "  +3.73%  "
$ws.Range("E12").Value = "  +11.17%  "
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("E17").Value = "  +6.06%  "
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("E23").Value = "  +5.63%  "
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("E27").Value = "  +2.94%  "
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("E29").Value = "  +3.23%  "
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("E31").Value = "  +4.36%  "
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("E33").Value = "  +6.73%  "
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  +3.51%  "
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("E39").Value = "  +4.41%  "
$ws.Range("E40").Value = "  +3.84%  "
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("E42").Value = "  +5.88%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("E47").Value = "  +4.54%  "
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("E50").Value = "  +4.99%  "
$ws.Range("E51").Value = "  +1.03%  "
